$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAR-22")

# --- Row 46: new entry (No.20, 2022-03-08) -------------------------------
# Copy formats from the row above (row 45) so borders/number-formats match,
# then overwrite with the new row's values.
$ws.Range("A45:F45").Copy() | Out-Null
$ws.Range("A46:F46").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A46").Value = 20
$ws.Range("B46").Value = 44648
$ws.Range("C46").Value = "RPA GSS"
$ws.Range("D46").WrapText = $true
$ws.Range("D46").Value = "A). Implementation of the public holidays is not the below mentioned tasks and instead of that the log information has been recorded `nat log file and they are tested and they running smoothly`n1. Credit_info`n2. Credit_details`n3.  DRS daily`n4. PR Summary daily`n5. Service Order cancelled"
$ws.Range("E46").Value = 1
$ws.Range("F46").Value = "Completed"
$ws.Rows.Item(46).RowHeight = 115.2

# --- Row 47: continuation line (RPA RLOGIC) ------------------------------
$ws.Range("C47").Value = "RPA RLOGIC"
$ws.Range("D47").Value = "1. MLR daily tasks are not executed due to password is matching and after receiving the passwords, we have rerun the daily task of MLR"
$ws.Range("E47").NumberFormat = "0%"
$ws.Range("E47").Value = 1
$ws.Range("F47").Value = "Completed"

# --- Row 48: continuation line (warranty / saw discount) -----------------
$ws.Range("D48").WrapText = $true
$ws.Range("D48").Value = "2. In addition to warranty daily task, in saw discount task also table structured  changed at GSPN and due to that no data found is `noccurred and it has been fixed by Kabilan san."
$ws.Range("E48").NumberFormat = "0%"
$ws.Range("E48").Value = 1
$ws.Range("F48").Value = "Completed"
$ws.Rows.Item(48).RowHeight = 28.8

# --- View state: match the author's final scroll/selection position -----
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$ws.Range("D52").Select() | Out-Null

Write-Output "done"
